$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text/value updates (dates, temperatures, pressures, radiation, wind, humidity) ---
$ws.Range('E2').Value = '2026-02-21 17:48:43'
$ws.Range('O2').Value = '4.1 °C'
$ws.Range('E3').Value = '2026-02-21 17:48:45'
$ws.Range('K3').Value = '16.2 MJ/m2'
$ws.Range('O3').Value = '1.8 °C'
$ws.Range('E4').Value = '2026-02-21 17:48:48'
$ws.Range('E5').Value = '2026-02-21 17:48:50'
$ws.Range('K5').Value = '16.1 MJ/m2'
$ws.Range('E6').Value = '2026-02-21 17:48:53'
$ws.Range('J6').Value = '1029.3 hPa'
$ws.Range('O6').Value = '10.5 °C'
$ws.Range('E7').Value = '2026-02-21 17:48:55'
$ws.Range('E8').Value = '2026-02-21 17:48:58'
$ws.Range('K8').Value = '16.1 MJ/m2'
$ws.Range('O8').Value = '10.8 °C'
$ws.Range('E9').Value = '2026-02-21 17:49:00'
$ws.Range('E10').Value = '2026-02-21 17:49:03'
$ws.Range('O10').Value = '8.8 °C'
$ws.Range('E11').Value = '2026-02-21 17:49:06'
$ws.Range('E12').Value = '2026-02-21 17:49:08'
$ws.Range('E13').Value = '2026-02-21 17:49:11'
$ws.Range('J13').Value = '1031.8 hPa'
$ws.Range('K13').Value = '16.1 MJ/m2'
$ws.Range('L13').Value = '18.0 km/h - 105º 17:26 TU'
$ws.Range('O13').Value = '5.1 °C'
$ws.Range('E14').Value = '2026-02-21 17:49:13'
$ws.Range('E15').Value = '2026-02-21 17:49:16'
$ws.Range('E16').Value = '2026-02-21 17:49:19'
$ws.Range('E17').Value = '2026-02-21 17:49:21'
$ws.Range('O17').Value = '8.6 °C'
$ws.Range('E18').Value = '2026-02-21 17:49:24'
$ws.Range('J18').Value = '1029.7 hPa'
$ws.Range('K18').Value = '15.4 MJ/m2'
$ws.Range('O18').Value = '8.7 °C'
$ws.Range('E19').Value = '2026-02-21 17:49:27'
$ws.Range('O19').Value = '7.9 °C'
$ws.Range('E20').Value = '2026-02-21 17:49:29'
$ws.Range('O20').Value = '3.1 °C'
$ws.Range('E21').Value = '2026-02-21 17:49:32'
$ws.Range('J21').Value = '1030.7 hPa'
$ws.Range('O21').Value = '7.3 °C'
$ws.Range('E22').Value = '2026-02-21 17:49:34'
$ws.Range('K22').Value = '16.7 MJ/m2'
$ws.Range('O22').Value = '1.8 °C'
$ws.Range('E23').Value = '2026-02-21 17:49:37'
$ws.Range('K23').Value = '16.1 MJ/m2'
$ws.Range('E24').Value = '2026-02-21 17:49:40'
$ws.Range('K24').Value = '15.9 MJ/m2'
$ws.Range('O24').Value = '6.3 °C'
$ws.Range('E25').Value = '2026-02-21 17:49:43'
$ws.Range('E26').Value = '2026-02-21 17:49:45'
$ws.Range('E27').Value = '2026-02-21 17:49:48'
$ws.Range('E28').Value = '2026-02-21 17:49:51'
$ws.Range('J28').Value = '1029.6 hPa'
$ws.Range('O28').Value = '8.2 °C'
$ws.Range('E29').Value = '2026-02-21 17:49:53'
$ws.Range('K29').Value = '15.4 MJ/m2'
$ws.Range('O29').Value = '12.2 °C'
$ws.Range('E30').Value = '2026-02-21 17:49:55'
$ws.Range('O30').Value = '12.1 °C'
$ws.Range('E31').Value = '2026-02-21 17:49:58'
$ws.Range('K31').Value = '15.2 MJ/m2'
$ws.Range('E32').Value = '2026-02-21 17:50:01'
$ws.Range('K32').Value = '16.2 MJ/m2'
$ws.Range('O32').Value = '5.9 °C'
$ws.Range('E33').Value = '2026-02-21 17:50:03'
$ws.Range('J33').Value = '1030.4 hPa'
$ws.Range('O33').Value = '6.3 °C'
$ws.Range('E34').Value = '2026-02-21 17:50:06'
$ws.Range('O34').Value = '5.1 °C'
$ws.Range('E35').Value = '2026-02-21 17:50:09'
$ws.Range('K35').Value = '16.5 MJ/m2'
$ws.Range('O35').Value = '8.0 °C'
$ws.Range('E36').Value = '2026-02-21 17:50:12'
$ws.Range('K36').Value = '15.3 MJ/m2'
$ws.Range('E37').Value = '2026-02-21 17:50:14'
$ws.Range('J37').Value = '1031.2 hPa'
$ws.Range('O37').Value = '6.0 °C'
$ws.Range('E38').Value = '2026-02-21 17:50:17'
$ws.Range('O38').Value = '9.8 °C'
$ws.Range('E39').Value = '2026-02-21 17:50:20'
$ws.Range('E40').Value = '2026-02-21 17:50:22'
$ws.Range('E41').Value = '2026-02-21 17:50:25'
$ws.Range('K41').Value = '15.5 MJ/m2'
$ws.Range('O41').Value = '11.4 °C'
$ws.Range('E42').Value = '2026-02-21 17:50:28'
$ws.Range('O42').Value = '11.1 °C'
$ws.Range('E43').Value = '2026-02-21 17:50:30'
$ws.Range('O43').Value = '6.7 °C'
$ws.Range('E44').Value = '2026-02-21 17:50:32'
$ws.Range('O44').Value = '2.6 °C'
$ws.Range('E45').Value = '2026-02-21 17:50:35'
$ws.Range('E46').Value = '2026-02-21 17:50:38'
$ws.Range('J46').Value = '1031.6 hPa'
$ws.Range('O46').Value = '10.0 °C'

# --- Percent-style text cells: Excel auto-converts plain "NN%" input into a
# numeric percentage, which would change both the stored value and the cell
# number format. To preserve these as literal text (matching the source data,
# which stores every column as plain text), prefix with an apostrophe to force
# text entry, then copy the number format back from an untouched same-style
# donor cell in the same column so the cell keeps its original "General" style.

$ws.Range('H18').Value = "'73%"
$ws.Range('H24').Value = "'82%"
$ws.Range('H32').Value = "'77%"
$ws.Range('H34').Value = "'36%"
$ws.Range('H38').Value = "'71%"
$ws.Range('H43').Value = "'77%"
$ws.Range('H44').Value = "'37%"
$ws.Range('H46').Value = "'67%"

$donor = $ws.Range('H2')
$donor.Copy()
$percentTargets = @('H18', 'H24', 'H32', 'H34', 'H38', 'H43', 'H44', 'H46')
foreach ($t in $percentTargets) {
  $ws.Range($t).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0
